$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 27, shifting existing rows 27:60 down to 28:61
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record
$ws.Cells.Item(27, 1).Value = 8
$ws.Cells.Item(27, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(27, 3).Value = "Coquimbo"
$ws.Cells.Item(27, 4).Value = 44803
$ws.Cells.Item(27, 5).Value = 4
$ws.Cells.Item(27, 6).Value = 100114007
$ws.Cells.Item(27, 7).Value = "Jengibre"
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 520
$ws.Cells.Item(27, 11).Value = 14000
$ws.Cells.Item(27, 12).Value = 15000
$ws.Cells.Item(27, 13).Value = 14500
$ws.Cells.Item(27, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(27, 15).Value = "Perú"
$ws.Cells.Item(27, 16).Value = 1115
$ws.Cells.Item(27, 17).Value = 13
$ws.Cells.Item(27, 18).Value = "Hortaliza"
